$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-coerced to a number by Excel
# (losing significant trailing zeros); force Text format before writing, then
# restore the default style so no stray formatting is left behind.
$textCells = @('D5', 'D6', 'D17', 'D19', 'D20', 'D21', 'D25', 'D31', 'D32', 'D33', 'D36', 'D37', 'D40', 'D41', 'D43', 'D47', 'D49', 'D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = '@'
}

$ws.Range('D2').Value = '68.361.93'
$ws.Range('E2').Value = '  -1.51%  '
$ws.Range('D3').Value = '2.447.40'
$ws.Range('E3').Value = '  -1.55%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '554.93'
$ws.Range('E5').Value = '  -2.46%  '
$ws.Range('D6').Value = '161.35'
$ws.Range('E6').Value = '  -1.75%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  -2.13%  '
$ws.Range('D9').Value = '2.446.01'
$ws.Range('E9').Value = '  -1.48%  '
$ws.Range('E10').Value = '  -7.15%  '
$ws.Range('E11').Value = '  -1.14%  '
$ws.Range('E12').Value = '  -5.84%  '
$ws.Range('E13').Value = '  -2.38%  '
$ws.Range('D14').Value = '2.896.63'
$ws.Range('E14').Value = '  -1.56%  '
$ws.Range('D15').Value = '68.241.72'
$ws.Range('E15').Value = '  -1.48%  '
$ws.Range('E16').Value = '  -4.43%  '
$ws.Range('D17').Value = '23.21'
$ws.Range('E17').Value = '  -3.96%  '
$ws.Range('D18').Value = '2.441.54'
$ws.Range('E18').Value = '  -2.18%  '
$ws.Range('D19').Value = '10.71'
$ws.Range('E19').Value = '  -3.66%  '
$ws.Range('D20').Value = '339.10'
$ws.Range('E20').Value = '  -2.08%  '
$ws.Range('D21').Value = '6.97'
$ws.Range('E21').Value = '  -5.31%  '
$ws.Range('E22').Value = '  -2.88%  '
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('E24').Value = '  -2.94%  '
$ws.Range('D25').Value = '66.26'
$ws.Range('E25').Value = '  -4.66%  '
$ws.Range('E26').Value = '  -6.43%  '
$ws.Range('D27').Value = '2.573.27'
$ws.Range('E27').Value = '  -1.84%  '
$ws.Range('E28').Value = '  +0.03%  '
$ws.Range('E29').Value = '  -6.66%  '
$ws.Range('D30').Value = '0.0₃0808'
$ws.Range('E30').Value = '  -6.59%  '
$ws.Range('D31').Value = '7.10'
$ws.Range('E31').Value = '  -6.44%  '
$ws.Range('B32').Value = 'Bittensor'
$ws.Range('C32').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D32').Value = '433.65'
$ws.Range('E32').Value = '  -0.62%  '
$ws.Range('B33').Value = 'FirstDigitalUSD'
$ws.Range('C33').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D33').Value = '0.999'
$ws.Range('E33').Value = '  +0.01%  '
$ws.Range('E34').Value = '  -5.53%  '
$ws.Range('E35').Value = '  -5.97%  '
$ws.Range('D36').Value = '157.42'
$ws.Range('E36').Value = '  +0.35%  '
$ws.Range('D37').Value = '19.01'
$ws.Range('E37').Value = '  -0.28%  '
$ws.Range('E38').Value = '  +0.14%  '
$ws.Range('E39').Value = '  -3.13%  '
$ws.Range('D40').Value = '17.74'
$ws.Range('E40').Value = '  -2.10%  '
$ws.Range('D41').Value = '0.300'
$ws.Range('E41').Value = '  -3.72%  '
$ws.Range('E42').Value = '  -3.52%  '
$ws.Range('D43').Value = '37.36'
$ws.Range('E43').Value = '  -0.93%  '
$ws.Range('E44').Value = '  -7.74%  '
$ws.Range('E45').Value = '  +1.16%  '
$ws.Range('E46').Value = '  -5.90%  '
$ws.Range('D47').Value = '131.20'
$ws.Range('E47').Value = '  -4.66%  '
$ws.Range('E48').Value = '  -2.83%  '
$ws.Range('D49').Value = '0.0709'
$ws.Range('E49').Value = '  -1.96%  '
$ws.Range('E50').Value = '  -5.94%  '
$ws.Range('D51').Value = '0.557'
$ws.Range('E51').Value = '  -2.40%  '

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = 'Normal'
}
